$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "SCHEMA FILE ASCII" column values for the Excel/Access/dBase/XML rows
$ws.Range("D2").Value = "none"
$ws.Range("D3").Value = "none"
$ws.Range("D4").Value = "none"
$ws.Range("D5").Value = "xsd , xrdf"

# New import-definition columns for Print Report & Adobe / Text / AS400
$ws.Range("C6").Value = "jpm"
$ws.Range("C8").Value = "fdf"
$ws.Range("C7").Value = "rdm"
$ws.Range("D7").Value = "rdf"
$ws.Range("D6").Value = "jpm but for ascii"
$ws.Range("D8").Value = "fdf"

# Move the active selection to F5 (import definition button placement)
$ws.Range("F5").Select()
